# Update LR-pair TPM-derived statistics on the first worksheet.
# This mirrors the author's "update scripts wuth new tpm" commit: the
# Ligand-expressing cell count for the Cntn6 -> Notch1 pair increased
# from 1 to 2, which cascades into the dependent detection-rate /
# expression / specificity columns for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07426566666666666
$ws.Range("H2").Value = 0.222797
$ws.Range("M2").Value = 38.55267666666666
$ws.Range("N2").Value = 115.65803
$ws.Range("O2").Value = 0.5758151725879548
$ws.Range("P2").Value = 0.5758151725879548
$ws.Range("Q2").Value = 2.863140234434444
$ws.Range("R2").Value = 25.76826210991
$ws.Range("S2").Value = 0.5758151725879548
$ws.Range("T2").Value = 0.5758151725879548

# --- Row 3 ---------------------------------------------------------
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07426566666666666
$ws.Range("H3").Value = 0.222797
$ws.Range("O3").Value = 0.08021535714867321
$ws.Range("P3").Value = 0.08021535714867323
$ws.Range("Q3").Value = 0.3988568335906666
$ws.Range("R3").Value = 3.589711502316
$ws.Range("S3").Value = 0.08021535714867321
$ws.Range("T3").Value = 0.08021535714867323

# --- Row 4 ---------------------------------------------------------
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07426566666666666
$ws.Range("H4").Value = 0.222797
$ws.Range("M4").Value = 23.02986166666667
$ws.Range("N4").Value = 69.089585
$ws.Range("O4").Value = 0.3439694702633719
$ws.Range("P4").Value = 0.3439694702633719
$ws.Range("Q4").Value = 1.710328029916111
$ws.Range("R4").Value = 15.392952269245
$ws.Range("S4").Value = 0.3439694702633719
$ws.Range("T4").Value = 0.3439694702633719
